$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the Weight_g (column C) values for the existing rows 362-370 ---
$ws.Cells.Item(362, 3).Value = 459
$ws.Cells.Item(363, 3).Value = 438
$ws.Cells.Item(364, 3).Value = 449
$ws.Cells.Item(365, 3).Value = 499
$ws.Cells.Item(366, 3).Value = 477
$ws.Cells.Item(367, 3).Value = 542
$ws.Cells.Item(368, 3).Value = 458
$ws.Cells.Item(369, 3).Value = 509
$ws.Cells.Item(370, 3).Value = 507

# --- Append four new day-blocks (rows 371-406) ---
# --- day block starting row 371 ---
$ws.Cells.Item(371, 1).Value = 1
$ws.Cells.Item(371, 2).Formula = "=B362+1"
$ws.Cells.Item(371, 3).Value = 467
$ws.Cells.Item(371, 4).Value = 1
$ws.Cells.Item(371, 5).Value = "Black tipped tail"
$ws.Cells.Item(371, 6).Value = "BTT"
$ws.Cells.Item(372, 1).Value = 2
$ws.Cells.Item(372, 2).Formula = "=B371"
$ws.Cells.Item(372, 3).Value = 440
$ws.Cells.Item(372, 4).Value = 1
$ws.Cells.Item(372, 5).Value = "Black dot left rear leg"
$ws.Cells.Item(372, 6).Value = "BDLRL"
$ws.Cells.Item(373, 1).Value = 3
$ws.Cells.Item(373, 2).Formula = "=B364+1"
$ws.Cells.Item(373, 3).Value = 449
$ws.Cells.Item(373, 4).Value = 1
$ws.Cells.Item(373, 5).Value = "White rear legs"
$ws.Cells.Item(373, 6).Value = "WRL"
$ws.Cells.Item(374, 1).Value = 4
$ws.Cells.Item(374, 2).Formula = "=B373"
$ws.Cells.Item(374, 3).Value = 499
$ws.Cells.Item(374, 4).Value = 2
$ws.Cells.Item(374, 5).Value = "Half white tail/Black spots under side"
$ws.Cells.Item(374, 6).Value = "1/2WT/BSU"
$ws.Cells.Item(375, 1).Value = 5
$ws.Cells.Item(375, 2).Formula = "=B366+1"
$ws.Cells.Item(375, 3).Value = 475
$ws.Cells.Item(375, 4).Value = 2
$ws.Cells.Item(375, 5).Value = "Half white tail/White under side"
$ws.Cells.Item(375, 6).Value = "1/2WT/WU"
$ws.Cells.Item(376, 1).Value = 6
$ws.Cells.Item(376, 2).Formula = "=B367+1"
$ws.Cells.Item(376, 3).Value = 548
$ws.Cells.Item(376, 4).Value = 2
$ws.Cells.Item(376, 5).Value = "Small white tip tail"
$ws.Cells.Item(376, 6).Value = "SWTT"
$ws.Cells.Item(377, 1).Value = 7
$ws.Cells.Item(377, 2).Formula = "=B376"
$ws.Cells.Item(377, 3).Value = 463
$ws.Cells.Item(377, 4).Value = 3
$ws.Cells.Item(377, 5).Value = "White tail"
$ws.Cells.Item(377, 6).Value = "WT"
$ws.Cells.Item(378, 1).Value = 8
$ws.Cells.Item(378, 2).Formula = "=B369+1"
$ws.Cells.Item(378, 3).Value = 515
$ws.Cells.Item(378, 4).Value = 3
$ws.Cells.Item(378, 5).Value = "Half white tail"
$ws.Cells.Item(378, 6).Value = "1/2WT"
$ws.Cells.Item(379, 1).Value = 9
$ws.Cells.Item(379, 2).Formula = "=B378"
$ws.Cells.Item(379, 3).Value = 513
$ws.Cells.Item(379, 4).Value = 3
$ws.Cells.Item(379, 5).Value = "One fourth white tip tail"
$ws.Cells.Item(379, 6).Value = "1/4WTT"
# --- day block starting row 380 ---
$ws.Cells.Item(380, 1).Value = 1
$ws.Cells.Item(380, 2).Formula = "=B371+1"
$ws.Cells.Item(380, 3).Value = 468
$ws.Cells.Item(380, 4).Value = 1
$ws.Cells.Item(380, 5).Value = "Black tipped tail"
$ws.Cells.Item(380, 6).Value = "BTT"
$ws.Cells.Item(381, 1).Value = 2
$ws.Cells.Item(381, 2).Formula = "=B380"
$ws.Cells.Item(381, 3).Value = 442
$ws.Cells.Item(381, 4).Value = 1
$ws.Cells.Item(381, 5).Value = "Black dot left rear leg"
$ws.Cells.Item(381, 6).Value = "BDLRL"
$ws.Cells.Item(382, 1).Value = 3
$ws.Cells.Item(382, 2).Formula = "=B373+1"
$ws.Cells.Item(382, 3).Value = 452
$ws.Cells.Item(382, 4).Value = 1
$ws.Cells.Item(382, 5).Value = "White rear legs"
$ws.Cells.Item(382, 6).Value = "WRL"
$ws.Cells.Item(383, 1).Value = 4
$ws.Cells.Item(383, 2).Formula = "=B382"
$ws.Cells.Item(383, 3).Value = 501
$ws.Cells.Item(383, 4).Value = 2
$ws.Cells.Item(383, 5).Value = "Half white tail/Black spots under side"
$ws.Cells.Item(383, 6).Value = "1/2WT/BSU"
$ws.Cells.Item(384, 1).Value = 5
$ws.Cells.Item(384, 2).Formula = "=B375+1"
$ws.Cells.Item(384, 3).Value = 481
$ws.Cells.Item(384, 4).Value = 2
$ws.Cells.Item(384, 5).Value = "Half white tail/White under side"
$ws.Cells.Item(384, 6).Value = "1/2WT/WU"
$ws.Cells.Item(385, 1).Value = 6
$ws.Cells.Item(385, 2).Formula = "=B376+1"
$ws.Cells.Item(385, 3).Value = 549
$ws.Cells.Item(385, 4).Value = 2
$ws.Cells.Item(385, 5).Value = "Small white tip tail"
$ws.Cells.Item(385, 6).Value = "SWTT"
$ws.Cells.Item(386, 1).Value = 7
$ws.Cells.Item(386, 2).Formula = "=B385"
$ws.Cells.Item(386, 3).Value = 470
$ws.Cells.Item(386, 4).Value = 3
$ws.Cells.Item(386, 5).Value = "White tail"
$ws.Cells.Item(386, 6).Value = "WT"
$ws.Cells.Item(387, 1).Value = 8
$ws.Cells.Item(387, 2).Formula = "=B378+1"
$ws.Cells.Item(387, 3).Value = 522
$ws.Cells.Item(387, 4).Value = 3
$ws.Cells.Item(387, 5).Value = "Half white tail"
$ws.Cells.Item(387, 6).Value = "1/2WT"
$ws.Cells.Item(388, 1).Value = 9
$ws.Cells.Item(388, 2).Formula = "=B387"
$ws.Cells.Item(388, 3).Value = 510
$ws.Cells.Item(388, 4).Value = 3
$ws.Cells.Item(388, 5).Value = "One fourth white tip tail"
$ws.Cells.Item(388, 6).Value = "1/4WTT"
# --- day block starting row 389 ---
$ws.Cells.Item(389, 1).Value = 1
$ws.Cells.Item(389, 2).Formula = "=B380+1"
$ws.Cells.Item(389, 3).Value = 470
$ws.Cells.Item(389, 4).Value = 1
$ws.Cells.Item(389, 5).Value = "Black tipped tail"
$ws.Cells.Item(389, 6).Value = "BTT"
$ws.Cells.Item(390, 1).Value = 2
$ws.Cells.Item(390, 2).Formula = "=B389"
$ws.Cells.Item(390, 3).Value = 444
$ws.Cells.Item(390, 4).Value = 1
$ws.Cells.Item(390, 5).Value = "Black dot left rear leg"
$ws.Cells.Item(390, 6).Value = "BDLRL"
$ws.Cells.Item(391, 1).Value = 3
$ws.Cells.Item(391, 2).Formula = "=B382+1"
$ws.Cells.Item(391, 3).Value = 454
$ws.Cells.Item(391, 4).Value = 1
$ws.Cells.Item(391, 5).Value = "White rear legs"
$ws.Cells.Item(391, 6).Value = "WRL"
$ws.Cells.Item(392, 1).Value = 4
$ws.Cells.Item(392, 2).Formula = "=B391"
$ws.Cells.Item(392, 3).Value = 500
$ws.Cells.Item(392, 4).Value = 2
$ws.Cells.Item(392, 5).Value = "Half white tail/Black spots under side"
$ws.Cells.Item(392, 6).Value = "1/2WT/BSU"
$ws.Cells.Item(393, 1).Value = 5
$ws.Cells.Item(393, 2).Formula = "=B384+1"
$ws.Cells.Item(393, 3).Value = 481
$ws.Cells.Item(393, 4).Value = 2
$ws.Cells.Item(393, 5).Value = "Half white tail/White under side"
$ws.Cells.Item(393, 6).Value = "1/2WT/WU"
$ws.Cells.Item(394, 1).Value = 6
$ws.Cells.Item(394, 2).Formula = "=B385+1"
$ws.Cells.Item(394, 3).Value = 553
$ws.Cells.Item(394, 4).Value = 2
$ws.Cells.Item(394, 5).Value = "Small white tip tail"
$ws.Cells.Item(394, 6).Value = "SWTT"
$ws.Cells.Item(395, 1).Value = 7
$ws.Cells.Item(395, 2).Formula = "=B394"
$ws.Cells.Item(395, 3).Value = 462
$ws.Cells.Item(395, 4).Value = 3
$ws.Cells.Item(395, 5).Value = "White tail"
$ws.Cells.Item(395, 6).Value = "WT"
$ws.Cells.Item(396, 1).Value = 8
$ws.Cells.Item(396, 2).Formula = "=B387+1"
$ws.Cells.Item(396, 3).Value = 518
$ws.Cells.Item(396, 4).Value = 3
$ws.Cells.Item(396, 5).Value = "Half white tail"
$ws.Cells.Item(396, 6).Value = "1/2WT"
$ws.Cells.Item(397, 1).Value = 9
$ws.Cells.Item(397, 2).Formula = "=B396"
$ws.Cells.Item(397, 3).Value = 513
$ws.Cells.Item(397, 4).Value = 3
$ws.Cells.Item(397, 5).Value = "One fourth white tip tail"
$ws.Cells.Item(397, 6).Value = "1/4WTT"
# --- day block starting row 398 ---
$ws.Cells.Item(398, 1).Value = 1
$ws.Cells.Item(398, 2).Formula = "=B389+1"
$ws.Cells.Item(398, 4).Value = 1
$ws.Cells.Item(398, 5).Value = "Black tipped tail"
$ws.Cells.Item(398, 6).Value = "BTT"
$ws.Cells.Item(399, 1).Value = 2
$ws.Cells.Item(399, 2).Formula = "=B398"
$ws.Cells.Item(399, 4).Value = 1
$ws.Cells.Item(399, 5).Value = "Black dot left rear leg"
$ws.Cells.Item(399, 6).Value = "BDLRL"
$ws.Cells.Item(400, 1).Value = 3
$ws.Cells.Item(400, 2).Formula = "=B391+1"
$ws.Cells.Item(400, 4).Value = 1
$ws.Cells.Item(400, 5).Value = "White rear legs"
$ws.Cells.Item(400, 6).Value = "WRL"
$ws.Cells.Item(401, 1).Value = 4
$ws.Cells.Item(401, 2).Formula = "=B400"
$ws.Cells.Item(401, 4).Value = 2
$ws.Cells.Item(401, 5).Value = "Half white tail/Black spots under side"
$ws.Cells.Item(401, 6).Value = "1/2WT/BSU"
$ws.Cells.Item(402, 1).Value = 5
$ws.Cells.Item(402, 2).Formula = "=B393+1"
$ws.Cells.Item(402, 4).Value = 2
$ws.Cells.Item(402, 5).Value = "Half white tail/White under side"
$ws.Cells.Item(402, 6).Value = "1/2WT/WU"
$ws.Cells.Item(403, 1).Value = 6
$ws.Cells.Item(403, 2).Formula = "=B394+1"
$ws.Cells.Item(403, 4).Value = 2
$ws.Cells.Item(403, 5).Value = "Small white tip tail"
$ws.Cells.Item(403, 6).Value = "SWTT"
$ws.Cells.Item(404, 1).Value = 7
$ws.Cells.Item(404, 2).Formula = "=B403"
$ws.Cells.Item(404, 4).Value = 3
$ws.Cells.Item(404, 5).Value = "White tail"
$ws.Cells.Item(404, 6).Value = "WT"
$ws.Cells.Item(405, 1).Value = 8
$ws.Cells.Item(405, 2).Formula = "=B396+1"
$ws.Cells.Item(405, 4).Value = 3
$ws.Cells.Item(405, 5).Value = "Half white tail"
$ws.Cells.Item(405, 6).Value = "1/2WT"
$ws.Cells.Item(406, 1).Value = 9
$ws.Cells.Item(406, 2).Formula = "=B405"
$ws.Cells.Item(406, 4).Value = 3
$ws.Cells.Item(406, 5).Value = "One fourth white tip tail"
$ws.Cells.Item(406, 6).Value = "1/4WTT"

# --- Update the view state to match where the author ended up editing ---
$ws.Range("C398").Select()
